# "restructured search results page for better processing"
#
# The sheet holds a City/State/Data lookup-results table (rows 2:52).
# Six rows are stale/duplicate entries that get dropped (everything below
# shifts up), and the new top-of-list rows get their "Data" status
# refreshed:
#   - row 2  (ARCTIC VLG / AL / not found)  -> removed entirely
#   - rows 12-14 (three duplicate Brooklyn / NY rows) -> removed
#   - row 27 (duplicate Chicago / IL row)    -> removed
#   - row 51 (duplicate Chicago / IL row)    -> removed
# After the shift, the new rows 2-4 (Chicago/IL, El Paso/TX, Los Angeles/CA)
# get their Data column populated: "added", "added", "not found".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid while we work.
$ws.Rows(51).Delete()
$ws.Rows(27).Delete()
$ws.Rows(14).Delete()
$ws.Rows(13).Delete()
$ws.Rows(12).Delete()
$ws.Rows(2).Delete()

# Refresh the "Data" column for the first three rows of what's left.
$ws.Range("C2").Value = "added"
$ws.Range("C3").Value = "added"
$ws.Range("C4").Value = "not found"
